$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Update task descriptions and effort (Time) values for the TanLoc estimate sheet.
$ws.Range("B2").Value = "Create New Product"
$ws.Range("C2").Value = 3.5

$ws.Range("B3").Value = "Design UI"

$ws.Range("B4").Value = "Write a function to loading all products"
$ws.Range("C4").Value = 1.5

$ws.Range("B5").Value = "Detailed coding  of product functions"
$ws.Range("C5").Value = 1.5

$ws.Range("B6").Value = "Unit Test"
$ws.Range("C6").Value = 1

# Remove the old row 7 (its task moved into row 6 / was superseded by "Unit Test");
# this shifts the "Sum" row up from row 8 to row 7 and keeps its SUM formulas
# correctly adjusted to the new data range (C2:C6 / D2:D6).
$ws.Rows.Item(7).Delete()
